$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 17-24, column E (Periodo Mora) and F (Valor Mora) are updated.
# The periods for rows 17-23 get reversed (2208..2202), and row 24 becomes 2201.
# The special amount 29333 (previously on period 2208 / row 24) now travels with
# period 2208 to row 17; all other rows revert to the standard 40000 value, except
# row 24 (period 2201) which now also uses 40000.

$ws.Range("E17").Value = "2208"
$ws.Range("F17").Value = 29333

$ws.Range("E18").Value = "2207"
$ws.Range("F18").Value = 40000

$ws.Range("E19").Value = "2206"
$ws.Range("F19").Value = 40000

$ws.Range("E20").Value = "2205"
$ws.Range("F20").Value = 40000

$ws.Range("E21").Value = "2204"
$ws.Range("F21").Value = 40000

$ws.Range("E22").Value = "2203"
$ws.Range("F22").Value = 40000

$ws.Range("E23").Value = "2202"
$ws.Range("F23").Value = 40000

$ws.Range("E24").Value = "2201"
$ws.Range("F24").Value = 40000
